$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0000599376648285783
$ws.Range("C2").Value = 0.999280748022057
$ws.Range("D2").Value = 0.101474466554783
$ws.Range("E2").Value = 0.974526492447854
$ws.Range("F2").Value = 0.0000599376648285783
$ws.Range("G2").Value = 0.0788180292495804
$ws.Range("H2").Value = 0.997003116758571
$ws.Range("I2").Value = 0.000419563653800048
$ws.Range("J2").Value = 0.976264684727883
$ws.Range("K2").Value = 0.00173819228002877
$ws.Range("L2").Value = 0.999640374011029
$ws.Range("M2").Value = 0.867597698393671
$ws.Range("N2").Value = 0.0000599376648285783
$ws.Range("P2").Value = 0.983217453847998
$ws.Range("Q2").Value = 0.9991608726924
$ws.Range("R2").Value = 0.000959002637257252
$ws.Range("S2").Value = 0.00743227043874371
$ws.Range("T2").Value = 0.000239750659314313
$ws.Range("U2").Value = 0.000299688324142891
$ws.Range("V2").Value = 0.991668664588828
$ws.Range("W2").Value = 0.999100935027571
$ws.Range("X2").Value = 0.000239750659314313
$ws.Range("B3").Value = 0.999820187005514
$ws.Range("C3").Value = 0.000179812994485735
$ws.Range("D3").Value = 0.0000599376648285783
$ws.Range("E3").Value = 0.0119275953008871
$ws.Range("F3").Value = 0.000179812994485735
$ws.Range("G3").Value = 0.000239750659314313
$ws.Range("H3").Value = 0.00215775593382882
$ws.Range("M3").Value = 0.0364421002157756
$ws.Range("N3").Value = 0.999640374011029
$ws.Range("P3").Value = 0.00179812994485735
$ws.Range("Q3").Value = 0.000419563653800048
$ws.Range("U3").Value = 0.000119875329657157
$ws.Range("V3").Value = 0.000239750659314313
$ws.Range("W3").Value = 0.0000599376648285783
$ws.Range("X3").Value = 0.000119875329657157
$ws.Range("C4").Value = 0.00035962598897147
$ws.Range("D4").Value = 0.895288899544474
$ws.Range("E4").Value = 0.00935027571325821
$ws.Range("F4").Value = 0.000239750659314313
$ws.Range("G4").Value = 0.918245025173819
$ws.Range("H4").Value = 0.000539438983457205
$ws.Range("I4").Value = 0.999040997362743
$ws.Range("J4").Value = 0.0224166866458883
$ws.Range("K4").Value = 0.998261807719971
$ws.Range("L4").Value = 0.000299688324142891
$ws.Range("M4").Value = 0.0837329177655239
$ws.Range("N4").Value = 0.0000599376648285783
$ws.Range("P4").Value = 0.00377607288420043
$ws.Range("Q4").Value = 0.000179812994485735
$ws.Range("R4").Value = 0.998981059697914
$ws.Range("S4").Value = 0.992148165907456
$ws.Range("T4").Value = 0.999760249340686
$ws.Range("U4").Value = 0.9995804363462
$ws.Range("V4").Value = 0.00719251977942939
$ws.Range("W4").Value = 0.000839127307600096
$ws.Range("X4").Value = 0.9995804363462
$ws.Range("B5").Value = 0.0000599376648285783
$ws.Range("C5").Value = 0.000179812994485735
$ws.Range("D5").Value = 0.0000599376648285783
$ws.Range("E5").Value = 0.00341644689522896
$ws.Range("F5").Value = 0.999400623351714
$ws.Range("G5").Value = 0.00035962598897147
$ws.Range("H5").Value = 0.000179812994485735
$ws.Range("I5").Value = 0.000299688324142891
$ws.Range("J5").Value = 0.0000599376648285783
$ws.Range("L5").Value = 0.0000599376648285783
$ws.Range("M5").Value = 0.00923040038360106
$ws.Range("N5").Value = 0.000239750659314313
$ws.Range("P5").Value = 0.0100095900263726
$ws.Range("Q5").Value = 0.000119875329657157
$ws.Range("R5").Value = 0.0000599376648285783
$ws.Range("V5").Value = 0.000179812994485735
$ws.Range("X5").Value = 0.0000599376648285783
